$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 7).Value = 21.18599966666667
$ws.Cells.Item(2, 8).Value = 63.557999
$ws.Cells.Item(2, 9).Value = 0.08765141600314529
$ws.Cells.Item(2, 10).Value = 0.08765141600314529
$ws.Cells.Item(2, 13).Value = 1.646992
$ws.Cells.Item(2, 14).Value = 4.940976
$ws.Cells.Item(2, 15).Value = 0.2071783517404009
$ws.Cells.Item(2, 16).Value = 0.2071783517404009
$ws.Cells.Item(2, 17).Value = 34.89317196300266
$ws.Cells.Item(2, 18).Value = 314.038547667024
$ws.Cells.Item(2, 19).Value = 0.01815947589524384
$ws.Cells.Item(2, 20).Value = 0.01815947589524384
$ws.Cells.Item(3, 7).Value = 21.18599966666667
$ws.Cells.Item(3, 8).Value = 63.557999
$ws.Cells.Item(3, 9).Value = 0.08765141600314529
$ws.Cells.Item(3, 10).Value = 0.08765141600314529
$ws.Cells.Item(3, 15).Value = 0.4685125322965616
$ws.Cells.Item(3, 16).Value = 0.4685125322965616
$ws.Cells.Item(3, 17).Value = 78.90731931649901
$ws.Cells.Item(3, 18).Value = 710.1658738484911
$ws.Cells.Item(3, 19).Value = 0.04106578687101296
$ws.Cells.Item(3, 20).Value = 0.04106578687101296
$ws.Cells.Item(4, 7).Value = 21.18599966666667
$ws.Cells.Item(4, 8).Value = 63.557999
$ws.Cells.Item(4, 9).Value = 0.08765141600314529
$ws.Cells.Item(4, 10).Value = 0.08765141600314529
$ws.Cells.Item(4, 13).Value = 2.284352333333333
$ws.Cells.Item(4, 14).Value = 6.853057
$ws.Cells.Item(4, 15).Value = 0.2873531572796583
$ws.Cells.Item(4, 16).Value = 0.2873531572796583
$ws.Cells.Item(4, 17).Value = 48.39628777254922
$ws.Cells.Item(4, 18).Value = 435.566589952943
$ws.Cells.Item(4, 19).Value = 0.02518691112853656
$ws.Cells.Item(4, 20).Value = 0.02518691112853656
$ws.Cells.Item(5, 7).Value = 21.18599966666667
$ws.Cells.Item(5, 8).Value = 63.557999
$ws.Cells.Item(5, 9).Value = 0.08765141600314529
$ws.Cells.Item(5, 10).Value = 0.08765141600314529
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.2937863333333333
$ws.Cells.Item(5, 14).Value = 0.881359
$ws.Cells.Item(5, 15).Value = 0.03695595868337916
$ws.Cells.Item(5, 16).Value = 0.03695595868337916
$ws.Cells.Item(5, 17).Value = 6.224157160071222
$ws.Cells.Item(5, 18).Value = 56.01741444064101
$ws.Cells.Item(5, 19).Value = 0.003239242108351916
$ws.Cells.Item(5, 20).Value = 0.003239242108351916
$ws.Cells.Item(6, 9).Value = 0.5040014103551328
$ws.Cells.Item(6, 10).Value = 0.5040014103551328
$ws.Cells.Item(6, 13).Value = 1.646992
$ws.Cells.Item(6, 14).Value = 4.940976
$ws.Cells.Item(6, 15).Value = 0.2071783517404009
$ws.Cells.Item(6, 16).Value = 0.2071783517404009
$ws.Cells.Item(6, 17).Value = 200.6380351058613
$ws.Cells.Item(6, 18).Value = 1805.742315952752
$ws.Cells.Item(6, 19).Value = 0.1044181814722139
$ws.Cells.Item(6, 20).Value = 0.1044181814722139
$ws.Cells.Item(7, 9).Value = 0.5040014103551328
$ws.Cells.Item(7, 10).Value = 0.5040014103551328
$ws.Cells.Item(7, 15).Value = 0.4685125322965616
$ws.Cells.Item(7, 16).Value = 0.4685125322965616
$ws.Cells.Item(7, 19).Value = 0.2361309770465218
$ws.Cells.Item(7, 20).Value = 0.2361309770465218
$ws.Cells.Item(8, 9).Value = 0.5040014103551328
$ws.Cells.Item(8, 10).Value = 0.5040014103551328
$ws.Cells.Item(8, 13).Value = 2.284352333333333
$ws.Cells.Item(8, 14).Value = 6.853057
$ws.Cells.Item(8, 15).Value = 0.2873531572796583
$ws.Cells.Item(8, 16).Value = 0.2873531572796583
$ws.Cells.Item(8, 17).Value = 278.2818396503988
$ws.Cells.Item(8, 18).Value = 2504.536556853589
$ws.Cells.Item(8, 19).Value = 0.144826396538948
$ws.Cells.Item(8, 20).Value = 0.144826396538948
$ws.Cells.Item(9, 9).Value = 0.5040014103551328
$ws.Cells.Item(9, 10).Value = 0.5040014103551328
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.2937863333333333
$ws.Cells.Item(9, 14).Value = 0.881359
$ws.Cells.Item(9, 15).Value = 0.03695595868337916
$ws.Cells.Item(9, 16).Value = 0.03695595868337916
$ws.Cells.Item(9, 17).Value = 35.78931328200478
$ws.Cells.Item(9, 18).Value = 322.103819538043
$ws.Cells.Item(9, 19).Value = 0.01862585529744911
$ws.Cells.Item(9, 20).Value = 0.01862585529744911
$ws.Cells.Item(10, 7).Value = 37.20718233333333
$ws.Cells.Item(10, 8).Value = 111.621547
$ws.Cells.Item(10, 9).Value = 0.1539347809079331
$ws.Cells.Item(10, 10).Value = 0.1539347809079331
$ws.Cells.Item(10, 13).Value = 1.646992
$ws.Cells.Item(10, 14).Value = 4.940976
$ws.Cells.Item(10, 15).Value = 0.2071783517404009
$ws.Cells.Item(10, 16).Value = 0.2071783517404009
$ws.Cells.Item(10, 17).Value = 61.27993164554132
$ws.Cells.Item(10, 18).Value = 551.519384809872
$ws.Cells.Item(10, 19).Value = 0.03189195418402532
$ws.Cells.Item(10, 20).Value = 0.03189195418402533
$ws.Cells.Item(11, 7).Value = 37.20718233333333
$ws.Cells.Item(11, 8).Value = 111.621547
$ws.Cells.Item(11, 9).Value = 0.1539347809079331
$ws.Cells.Item(11, 10).Value = 0.1539347809079331
$ws.Cells.Item(11, 15).Value = 0.4685125322965616
$ws.Cells.Item(11, 16).Value = 0.4685125322965616
$ws.Cells.Item(11, 17).Value = 138.578262222047
$ws.Cells.Item(11, 18).Value = 1247.204359998423
$ws.Cells.Item(11, 19).Value = 0.07212037401169215
$ws.Cells.Item(11, 20).Value = 0.07212037401169216
$ws.Cells.Item(12, 7).Value = 37.20718233333333
$ws.Cells.Item(12, 8).Value = 111.621547
$ws.Cells.Item(12, 9).Value = 0.1539347809079331
$ws.Cells.Item(12, 10).Value = 0.1539347809079331
$ws.Cells.Item(12, 13).Value = 2.284352333333333
$ws.Cells.Item(12, 14).Value = 6.853057
$ws.Cells.Item(12, 15).Value = 0.2873531572796583
$ws.Cells.Item(12, 16).Value = 0.2873531572796583
$ws.Cells.Item(12, 17).Value = 84.99431377990875
$ws.Cells.Item(12, 18).Value = 764.9488240191789
$ws.Cells.Item(12, 19).Value = 0.04423364530904704
$ws.Cells.Item(12, 20).Value = 0.04423364530904705
$ws.Cells.Item(13, 7).Value = 37.20718233333333
$ws.Cells.Item(13, 8).Value = 111.621547
$ws.Cells.Item(13, 9).Value = 0.1539347809079331
$ws.Cells.Item(13, 10).Value = 0.1539347809079331
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.2937863333333333
$ws.Cells.Item(13, 14).Value = 0.881359
$ws.Cells.Item(13, 15).Value = 0.03695595868337916
$ws.Cells.Item(13, 16).Value = 0.03695595868337916
$ws.Cells.Item(13, 17).Value = 10.93096167137478
$ws.Cells.Item(13, 18).Value = 98.378655042373
$ws.Cells.Item(13, 19).Value = 0.005688807403168599
$ws.Cells.Item(13, 20).Value = 0.005688807403168601
$ws.Cells.Item(14, 7).Value = 61.49336899999999
$ws.Cells.Item(14, 8).Value = 184.480107
$ws.Cells.Item(14, 9).Value = 0.2544123927337887
$ws.Cells.Item(14, 10).Value = 0.2544123927337887
$ws.Cells.Item(14, 13).Value = 1.646992
$ws.Cells.Item(14, 14).Value = 4.940976
$ws.Cells.Item(14, 15).Value = 0.2071783517404009
$ws.Cells.Item(14, 16).Value = 0.2071783517404009
$ws.Cells.Item(14, 17).Value = 101.279086796048
$ws.Cells.Item(14, 18).Value = 911.5117811644319
$ws.Cells.Item(14, 19).Value = 0.05270874018891791
$ws.Cells.Item(14, 20).Value = 0.05270874018891791
$ws.Cells.Item(15, 7).Value = 61.49336899999999
$ws.Cells.Item(15, 8).Value = 184.480107
$ws.Cells.Item(15, 9).Value = 0.2544123927337887
$ws.Cells.Item(15, 10).Value = 0.2544123927337887
$ws.Cells.Item(15, 15).Value = 0.4685125322965616
$ws.Cells.Item(15, 16).Value = 0.4685125322965616
$ws.Cells.Item(15, 17).Value = 229.032237320607
$ws.Cells.Item(15, 18).Value = 2061.290135885463
$ws.Cells.Item(15, 19).Value = 0.1191953943673347
$ws.Cells.Item(15, 20).Value = 0.1191953943673347
$ws.Cells.Item(16, 7).Value = 61.49336899999999
$ws.Cells.Item(16, 8).Value = 184.480107
$ws.Cells.Item(16, 9).Value = 0.2544123927337887
$ws.Cells.Item(16, 10).Value = 0.2544123927337887
$ws.Cells.Item(16, 13).Value = 2.284352333333333
$ws.Cells.Item(16, 14).Value = 6.853057
$ws.Cells.Item(16, 15).Value = 0.2873531572796583
$ws.Cells.Item(16, 16).Value = 0.2873531572796583
$ws.Cells.Item(16, 17).Value = 140.4725209596776
$ws.Cells.Item(16, 18).Value = 1264.252688637099
$ws.Cells.Item(16, 19).Value = 0.07310620430312659
$ws.Cells.Item(16, 20).Value = 0.07310620430312659
$ws.Cells.Item(17, 7).Value = 61.49336899999999
$ws.Cells.Item(17, 8).Value = 184.480107
$ws.Cells.Item(17, 9).Value = 0.2544123927337887
$ws.Cells.Item(17, 10).Value = 0.2544123927337887
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.2937863333333333
$ws.Cells.Item(17, 14).Value = 0.881359
$ws.Cells.Item(17, 15).Value = 0.03695595868337916
$ws.Cells.Item(17, 16).Value = 0.03695595868337916
$ws.Cells.Item(17, 17).Value = 18.06591140282366
$ws.Cells.Item(17, 18).Value = 162.593202625413
$ws.Cells.Item(17, 19).Value = 0.009402053874409528
$ws.Cells.Item(17, 20).Value = 0.00940205387440953
